# Insert a new data row at row 65 (pushes existing rows 65-99 down to 66-100)
# then populate the new row with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 65:99 down by one row, duplicating row 65's formatting for the new row.
$ws.Rows("65:65").Insert()

# Populate the newly inserted row 65 with the new record.
$ws.Range("A65").Value = 4
$ws.Range("B65").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C65").Value = "Los Lagos"
$ws.Range("D65").Value = 44460
$ws.Range("D65").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E65").Value = 10
$ws.Range("F65").Value = "Fruta"
$ws.Range("G65").Value = 100102
$ws.Range("H65").Value = "Cítricos"
$ws.Range("I65").Value = 100102004
$ws.Range("J65").Value = "Mandarina"
$ws.Range("K65").Value = "Murcott"
$ws.Range("L65").Value = "Primera"
$ws.Range("M65").Value = 400
$ws.Range("N65").Value = 6000
$ws.Range("O65").Value = 6000
$ws.Range("P65").Value = 6000
$ws.Range("Q65").Value = "`$/bandeja 10 kilos"
$ws.Range("R65").Value = "Provincia de Limarí"
$ws.Range("S65").Value = 600
$ws.Range("T65").Value = 10
